$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L17").Value = 4172.3079
$ws.Range("J17").Value = 1390.7693
$ws.Range("N17").Value = -4508.3079
$ws.Range("H17").Value = 1308.9286
$ws.Range("H88").Value = 6925.2
$ws.Range("M88").Value = -5044.75
$ws.Range("I88").Value = 5450.75
$ws.Range("K88").Value = 5450.75
$ws.Range("I91").Value = 5450.75
$ws.Range("K91").Value = 5450.75
$ws.Range("M91").Value = -4046.75
$ws.Range("H91").Value = 6925.2
$ws.Range("N113").Value = -9795.375
$ws.Range("J113").Value = 3287.375
$ws.Range("H113").Value = 3079.9
$ws.Range("L113").Value = 3287.375
$ws.Range("M127").Value = 2691.5263
$ws.Range("K127").Value = 2268.4737
$ws.Range("H127").Value = 756.1579
$ws.Range("I127").Value = 756.1579
$ws.Range("M137").Value = -93750312
$ws.Range("H137").Value = 19240774
$ws.Range("K137").Value = 93752862
$ws.Range("L137").Value = 73452
$ws.Range("J137").Value = 24484
$ws.Range("N137").Value = -78552
$ws.Range("I137").Value = 31250954
$ws.Range("H138").Value = 3090.1428
$ws.Range("I138").Value = 2220.8572
$ws.Range("M138").Value = -1522.571599999999
$ws.Range("K138").Value = 6662.571599999999
$ws.Range("N141").Value = -18610
$ws.Range("H141").Value = 2033.8
$ws.Range("I141").Value = 1556.3334
$ws.Range("J141").Value = 2750
$ws.Range("K141").Value = 4669.0002
$ws.Range("L141").Value = 8250
$ws.Range("M141").Value = 510.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 59726.105
$ws.Range("K32").Value = 66399.7
$ws.Range("M32").Value = -66112.7
$ws.Range("I32").Value = 66399.7
$ws.Range("J80").Value = 151723
$ws.Range("L80").Value = 151723
$ws.Range("H80").Value = 151723
$ws.Range("N80").Value = -153719
$ws.Range("J83").Value = 151723
$ws.Range("H83").Value = 151723
$ws.Range("L83").Value = 455169
$ws.Range("N83").Value = -465153

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 185885.73
$ws.Range("J86").Value = 402457.4
$ws.Range("L86").Value = 402457.4
$ws.Range("N86").Value = -404703.4
$ws.Range("J89").Value = 402457.4
$ws.Range("H89").Value = 185885.73
$ws.Range("L89").Value = 2012287
$ws.Range("N89").Value = -2023519
$ws.Range("H105").Value = 43490444
$ws.Range("M105").Value = -47630509
$ws.Range("I105").Value = 47632256
$ws.Range("K105").Value = 47632256
$ws.Range("L111").Value = 99995
$ws.Range("J111").Value = 99995
$ws.Range("H111").Value = 99995
$ws.Range("N111").Value = -108175
$ws.Range("I134").Value = 1802.5
$ws.Range("K134").Value = 5407.5
$ws.Range("M134").Value = -2872.5
$ws.Range("H134").Value = 4128.391

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -125001345
$ws.Range("I31").Value = 125001640
$ws.Range("K31").Value = 125001640
$ws.Range("H31").Value = 35718176
$ws.Range("M34").Value = -125001438
$ws.Range("K34").Value = 125001640
$ws.Range("H34").Value = 35718176
$ws.Range("I34").Value = 125001640
$ws.Range("K38").Value = 15087.5
$ws.Range("I38").Value = 15087.5
$ws.Range("J38").Value = 18686.666
$ws.Range("M38").Value = -14710.5
$ws.Range("L38").Value = 18686.666
$ws.Range("H38").Value = 16630
$ws.Range("N38").Value = -19440.666
$ws.Range("N46").Value = -19108.666
$ws.Range("L46").Value = 18686.666
$ws.Range("K46").Value = 15087.5
$ws.Range("M46").Value = -14876.5
$ws.Range("J46").Value = 18686.666
$ws.Range("I46").Value = 15087.5
$ws.Range("H46").Value = 16630
$ws.Range("I58").Value = 3161
$ws.Range("K58").Value = 3161
$ws.Range("J58").Value = 5880.909
$ws.Range("L58").Value = 5880.909
$ws.Range("N58").Value = -6286.909
$ws.Range("M58").Value = -2958
$ws.Range("H58").Value = 4520.9546
$ws.Range("H86").Value = 12250
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("K86").Value = 12250
$ws.Range("M86").Value = -11127
$ws.Range("I86").Value = 12250
$ws.Range("N86").ClearContents()
$ws.Range("M89").Value = -55634
$ws.Range("J89").Value = 0
$ws.Range("H89").Value = 12250
$ws.Range("I89").Value = 12250
$ws.Range("L89").Value = 0
$ws.Range("K89").Value = 61250
$ws.Range("N89").ClearContents()
$ws.Range("I94").Value = 3847.3572
$ws.Range("N94").Value = -1581.375
$ws.Range("M94").Value = -3396.3572
$ws.Range("J94").Value = 679.375
$ws.Range("H94").Value = 2695.3635
$ws.Range("L94").Value = 679.375
$ws.Range("K94").Value = 3847.3572
$ws.Range("I134").Value = 12340.25
$ws.Range("K134").Value = 37020.75
$ws.Range("M134").Value = -34485.75
$ws.Range("H134").Value = 10401.526
$ws.Range("J136").Value = 5880.909
$ws.Range("N136").Value = -22742.727
$ws.Range("H136").Value = 4520.9546
$ws.Range("L136").Value = 17642.727
$ws.Range("K136").Value = 9483
$ws.Range("I136").Value = 3161
$ws.Range("M136").Value = -6933

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 860.1111
$ws.Range("N50").Value = -9962
$ws.Range("L50").Value = 9000
$ws.Range("J50").Value = 3000
$ws.Range("J53").Value = 3000
$ws.Range("N53").Value = -9962
$ws.Range("L53").Value = 9000
$ws.Range("H53").Value = 860.1111
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -5441
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("L74").Value = 28500
$ws.Range("H74").Value = 9256.5
$ws.Range("J74").Value = 9500
$ws.Range("N74").Value = -30622
$ws.Range("H77").Value = 9256.5
$ws.Range("N77").Value = -96108
$ws.Range("L77").Value = 85500
$ws.Range("J77").Value = 9500
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 15000
$ws.Range("H80").Value = 5000
$ws.Range("N80").Value = -16872
$ws.Range("J83").Value = 5000
$ws.Range("H83").Value = 5000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54360
$ws.Range("K128").Value = 1309886.52
$ws.Range("I128").Value = 436628.84
$ws.Range("H128").Value = 436628.84
$ws.Range("M128").Value = -1304906.52
$ws.Range("J132").Value = 2998.5
$ws.Range("N132").Value = -32046.5
$ws.Range("L132").Value = 26986.5
$ws.Range("I132").Value = 1399.5
$ws.Range("H132").Value = 1666
$ws.Range("M132").Value = -10065.5
$ws.Range("K132").Value = 12595.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L57").Value = 69995
$ws.Range("N57").Value = -71635
$ws.Range("J57").Value = 69995
$ws.Range("H57").Value = 69995
$ws.Range("J70").Value = 26665.666
$ws.Range("I70").Value = 45000
$ws.Range("N70").Value = -27205.666
$ws.Range("H70").Value = 31249.25
$ws.Range("K70").Value = 45000
$ws.Range("M70").Value = -44730
$ws.Range("L70").Value = 26665.666
$ws.Range("H73").Value = 31249.25
$ws.Range("K73").Value = 45000
$ws.Range("N73").Value = -28537.666
$ws.Range("I73").Value = 45000
$ws.Range("M73").Value = -44064
$ws.Range("J73").Value = 26665.666
$ws.Range("L73").Value = 26665.666
$ws.Range("I132").Value = 2601.5
$ws.Range("H132").Value = 5062.3335
$ws.Range("M132").Value = -5274.5
$ws.Range("K132").Value = 7804.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10671.556
$ws.Range("N7").Value = -15474.125
$ws.Range("J7").Value = 15250.125
$ws.Range("L7").Value = 15250.125
$ws.Range("N16").Value = -5153
$ws.Range("H16").Value = 2582.5
$ws.Range("J16").Value = 4813
$ws.Range("L16").Value = 4813
$ws.Range("I16").Value = 1467.25
$ws.Range("M16").Value = -1297.25
$ws.Range("K16").Value = 1467.25
$ws.Range("H40").Value = 10681.6
$ws.Range("I40").Value = 9850.591
$ws.Range("M40").Value = -9714.591
$ws.Range("K40").Value = 9850.591
$ws.Range("H126").Value = 10671.556
$ws.Range("J126").Value = 15250.125
$ws.Range("L126").Value = 45750.375
$ws.Range("N126").Value = -50690.375
$ws.Range("J136").Value = 5067.2
$ws.Range("N136").Value = -20301.6
$ws.Range("H136").Value = 4834.9473
$ws.Range("L136").Value = 15201.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K52").Value = 4285.2856
$ws.Range("H52").Value = 10979.6
$ws.Range("I52").Value = 4285.2856
$ws.Range("M52").Value = -4059.2856
$ws.Range("N140").Value = -101444.8
$ws.Range("L140").Value = 91084.8
$ws.Range("J140").Value = 91084.8
$ws.Range("H140").Value = 91084.8
